$d = $word.ActiveDocument

$replacements = @(
    @{old="786÷5="; new="735÷3="},
    @{old="920÷7="; new="991÷2="},
    @{old="468÷3="; new="467÷4="},
    @{old="596÷4="; new="494÷5="},
    @{old="519÷8="; new="115÷7="},
    @{old="647÷6="; new="648÷7="},
    @{old="824÷6="; new="191÷2="},
    @{old="449÷9="; new="727÷7="},
    @{old="271÷8="; new="959÷3="},
    @{old="170÷9="; new="238÷8="},
    @{old="580÷3="; new="116÷4="},
    @{old="898÷7="; new="882÷4="},
    @{old="172÷8="; new="532÷4="},
    @{old="276÷9="; new="493÷4="},
    @{old="856÷4="; new="507÷6="},
    @{old="357÷6="; new="531÷7="},
    @{old="134÷6="; new="272÷2="},
    @{old="808÷2="; new="324÷8="},
    @{old="474÷6="; new="844÷9="},
    @{old="858÷7="; new="753÷5="},
    @{old="400÷9="; new="763÷7="},
    @{old="765÷8="; new="703÷7="},
    @{old="527÷4="; new="890÷6="},
    @{old="905÷4="; new="568÷9="},
    @{old="488÷7="; new="288÷7="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $r.new, 2)
}
